# "Add choose input table": the existing 6-column table (with a label
# column A + data columns B:H, where column H just repeated column G)
# becomes a plain 7-column table: the corner cell A1 is left blank (the
# "D0" label is removed) and the redundant trailing column H is dropped.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# A1 held the "D0" shared string; it becomes an empty corner cell.
$ws.Range("A1").ClearContents()

# Column H duplicated column G's values/header and is removed entirely,
# shifting the used range from A1:H6 down to A1:G6.
$ws.Columns("H").Delete()

# Selection moves from G8 to B1.
$ws.Range("B1").Select()
